# This sheet is a weekly "Apio" (celery) price log for Femacal de La Calera.
# The commit adds one new (more recent) weekly record at the top of the
# data block (row 293) and pushes every existing record down by one row,
# with the oldest record (which fell off the bottom) landing in a brand
# new last row (357).
#
# Columns A,B,C,E,F,G,H,Q,R are constant for every data row in this sheet,
# so only D (Fecha) and I..P (Calidad..Precio $/Kg) need to move. We do
# this as a block shift using Range.Value2 so we don't have to hand-type
# every one of the ~64 rows involved, then patch in the brand-new record
# for row 293 and the constant columns for the newly created row 357.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 293
$lastDataRow  = 356
$newLastRow   = 357

# Snapshot the shifting columns before they get overwritten.
$dateCol   = $ws.Range("D$firstDataRow`:D$lastDataRow").Value2
$restCols  = $ws.Range("I$firstDataRow`:P$lastDataRow").Value2

# Shift everything down by one data row.
$ws.Range("D$($firstDataRow+1)`:D$newLastRow").Value2 = $dateCol
$ws.Range("I$($firstDataRow+1)`:P$newLastRow").Value2 = $restCols

# The shifted-in D357 cell needs the same date number format as the rest
# of the Fecha column (the block Value2 write above doesn't carry styles).
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D2").NumberFormat

# Row 293 becomes the brand-new record reported in this edit.
$ws.Cells.Item($firstDataRow, 4).Value2  = 44637     # Fecha
$ws.Cells.Item($firstDataRow, 9).Value2  = "Primera" # Calidad
$ws.Cells.Item($firstDataRow, 10).Value2 = 310        # Volumen
$ws.Cells.Item($firstDataRow, 11).Value2 = 9500       # Precio minimo
$ws.Cells.Item($firstDataRow, 12).Value2 = 10000      # Precio maximo
$ws.Cells.Item($firstDataRow, 13).Value2 = 9758       # Precio promedio ponderado
$ws.Cells.Item($firstDataRow, 16).Value2 = 1626       # Precio $/Kg

# Row 357 is a brand new row, so the columns that are constant across the
# whole sheet need to be written explicitly too.
$ws.Cells.Item($newLastRow, 1).Value2  = 3
$ws.Cells.Item($newLastRow, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item($newLastRow, 3).Value2  = "Coquimbo"
$ws.Cells.Item($newLastRow, 5).Value2  = 5
$ws.Cells.Item($newLastRow, 6).Value2  = 100112017
$ws.Cells.Item($newLastRow, 7).Value2  = "Apio"
$ws.Cells.Item($newLastRow, 8).Value2  = "Americana (o)"
$ws.Cells.Item($newLastRow, 17).Value2 = 6
$ws.Cells.Item($newLastRow, 18).Value2 = "Hortaliza"

Write-Output "Shift complete: rows $firstDataRow..$lastDataRow moved to $($firstDataRow+1)..$newLastRow; new record written at row $firstDataRow."
